$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value2 = 0.1867219917012448
$ws.Range("C2").Value2 = 0.5726141078838174
$ws.Range("J2").Value2 = 0.02074688796680498
$ws.Range("P2").Value2 = 0.1369294605809129
$ws.Range("S2").Value2 = 0.08298755186721991
$ws.Range("B3").Value2 = 0.007142857142857143
$ws.Range("C3").Value2 = 0.03571428571428571
$ws.Range("J3").Value2 = 0.06428571428571428
$ws.Range("P3").Value2 = 0.7142857142857143
$ws.Range("S3").Value2 = 0.1785714285714286
$ws.Range("J4").Value2 = 0.05555555555555555
$ws.Range("P4").Value2 = 0.6944444444444444
$ws.Range("S4").Value2 = 0.25
$ws.Range("B6").Value2 = 0.03669724770642202
$ws.Range("D6").Value2 = 0.009174311926605505
$ws.Range("F6").Value2 = 0.04128440366972477
$ws.Range("J6").Value2 = 0.2568807339449541
$ws.Range("O6").Value2 = 0.02293577981651376
$ws.Range("Q6").Value2 = 0.1605504587155963
$ws.Range("R6").Value2 = 0.08256880733944955
$ws.Range("S6").Value2 = 0.3899082568807339
$ws.Range("B7").Value2 = 0.105
$ws.Range("D7").Value2 = 0.02
$ws.Range("F7").Value2 = 0.045
$ws.Range("J7").Value2 = 0.12
$ws.Range("O7").Value2 = 0.025
$ws.Range("Q7").Value2 = 0.185
$ws.Range("R7").Value2 = 0.065
$ws.Range("S7").Value2 = 0.435
$ws.Range("B8").Value2 = 0.07526881720430108
$ws.Range("D8").Value2 = 0.01505376344086022
$ws.Range("F8").Value2 = 0.08602150537634409
$ws.Range("J8").Value2 = 0.0989247311827957
$ws.Range("O8").Value2 = 0.01290322580645161
$ws.Range("Q8").Value2 = 0.189247311827957
$ws.Range("R8").Value2 = 0.1268817204301075
$ws.Range("S8").Value2 = 0.3956989247311828
$ws.Range("B9").Value2 = 0.08227848101265822
$ws.Range("D9").Value2 = 0.01265822784810127
$ws.Range("E9").Value2 = 0.006329113924050633
$ws.Range("F9").Value2 = 0.06329113924050633
$ws.Range("J9").Value2 = 0.120253164556962
$ws.Range("O9").Value2 = 0.01265822784810127
$ws.Range("Q9").Value2 = 0.1835443037974684
$ws.Range("R9").Value2 = 0.120253164556962
$ws.Range("S9").Value2 = 0.3987341772151899
$ws.Range("B10").Value2 = 0.1020583190394511
$ws.Range("D10").Value2 = 0.01886792452830189
$ws.Range("E10").Value2 = 0.001715265866209262
$ws.Range("F10").Value2 = 0.06946826758147513
$ws.Range("J10").Value2 = 0.1157804459691252
$ws.Range("O10").Value2 = 0.02229845626072041
$ws.Range("Q10").Value2 = 0.1963979416809606
$ws.Range("R10").Value2 = 0.09348198970840481
$ws.Range("S10").Value2 = 0.3799313893653516
$ws.Range("G11").Value2 = 0.1644295302013423
$ws.Range("J11").Value2 = 0.04697986577181208
$ws.Range("K11").Value2 = 0.2114093959731544
$ws.Range("L11").Value2 = 0.5604026845637584
$ws.Range("S11").Value2 = 0.01677852348993289
$ws.Range("G12").Value2 = 0.7396449704142012
$ws.Range("J12").Value2 = 0.1952662721893491
$ws.Range("L12").Value2 = 0.01183431952662722
$ws.Range("S12").Value2 = 0.05325443786982249
$ws.Range("G13").Value2 = 0.7142857142857143
$ws.Range("J13").Value2 = 0.2142857142857143
$ws.Range("S13").Value2 = 0.07142857142857142
$ws.Range("F15").Value2 = 0.01492537313432836
$ws.Range("H15").Value2 = 0.1890547263681592
$ws.Range("I15").Value2 = 0.03482587064676617
$ws.Range("J15").Value2 = 0.4029850746268657
$ws.Range("K15").Value2 = 0.05472636815920398
$ws.Range("M15").Value2 = 0.009950248756218905
$ws.Range("O15").Value2 = 0.04477611940298507
$ws.Range("S15").Value2 = 0.2487562189054726
$ws.Range("F16").Value2 = 0.03289473684210526
$ws.Range("H16").Value2 = 0.1644736842105263
$ws.Range("I16").Value2 = 0.09868421052631579
$ws.Range("J16").Value2 = 0.3684210526315789
$ws.Range("K16").Value2 = 0.1644736842105263
$ws.Range("M16").Value2 = 0.02631578947368421
$ws.Range("O16").Value2 = 0.03947368421052631
$ws.Range("S16").Value2 = 0.1052631578947368
$ws.Range("F17").Value2 = 0.01442307692307692
$ws.Range("H17").Value2 = 0.2043269230769231
$ws.Range("I17").Value2 = 0.09375
$ws.Range("J17").Value2 = 0.4230769230769231
$ws.Range("K17").Value2 = 0.1057692307692308
$ws.Range("M17").Value2 = 0.01442307692307692
$ws.Range("O17").Value2 = 0.05528846153846154
$ws.Range("S17").Value2 = 0.0889423076923077
$ws.Range("F18").Value2 = 0.0228310502283105
$ws.Range("H18").Value2 = 0.1963470319634703
$ws.Range("I18").Value2 = 0.0821917808219178
$ws.Range("J18").Value2 = 0.4474885844748858
$ws.Range("K18").Value2 = 0.0867579908675799
$ws.Range("O18").Value2 = 0.0502283105022831
$ws.Range("S18").Value2 = 0.1141552511415525
$ws.Range("F19").Value2 = 0.0199501246882793
$ws.Range("H19").Value2 = 0.229426433915212
$ws.Range("I19").Value2 = 0.0714879467996675
$ws.Range("J19").Value2 = 0.3466334164588529
$ws.Range("K19").Value2 = 0.1097256857855362
$ws.Range("M19").Value2 = 0.02660016625103907
$ws.Range("N19").Value2 = 0.0008312551953449709
$ws.Range("O19").Value2 = 0.0714879467996675
$ws.Range("S19").Value2 = 0.1288445552784705
